$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 10
    4  = 3
    5  = 5
    6  = 4
    7  = 5
    8  = 4
    9  = 3
    10 = 6
    11 = 3
    12 = 4
    13 = 3
    14 = 2
    15 = 4
    16 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
